# Fix sorting and generate viable xlsx and charts
# Update the Avg_Time_ms values for the first two rows of the Data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("D2").Value = 0.67730754
$ws.Range("D3").Value = 1.488243
